$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3723.5833
$ws.Range("I33").Value = 4566.5557
$ws.Range("J33").Value = 1194.6666
$ws.Range("K33").Value = 4566.5557
$ws.Range("L33").Value = 1194.6666
$ws.Range("M33").Value = -4337.5557
$ws.Range("N33").Value = -1652.6666

$ws.Range("H43").Value = 1135
$ws.Range("I43").Value = 1127.5
$ws.Range("J43").Value = 1150
$ws.Range("K43").Value = 1127.5
$ws.Range("L43").Value = 1150
$ws.Range("M43").Value = -1058.5
$ws.Range("N43").Value = -1288

$ws.Range("H103").Value = 551.35297
$ws.Range("I103").Value = 656.2727
$ws.Range("J103").Value = 359
$ws.Range("K103").Value = 1968.8181
$ws.Range("L103").Value = 1077
$ws.Range("M103").Value = -1382.8181
$ws.Range("N103").Value = -2249

$ws.Range("H110").Value = 13966.667
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws.Range("H111").Value = 696.8
$ws.Range("I111").Value = 621
$ws.Range("K111").Value = 1863
$ws.Range("M111").Value = 1204

$ws.Range("H115").Value = 470
$ws.Range("I115").Value = 455
$ws.Range("K115").Value = 1365
$ws.Range("M115").Value = 202

$ws.Range("H135").Value = 1790.5834
$ws.Range("I135").Value = 1720.8889
$ws.Range("K135").Value = 15488.0001
$ws.Range("M135").Value = -12953.0001

$ws.Range("H138").Value = 3119.1482
$ws.Range("I138").Value = 1750.3684
$ws.Range("J138").Value = 3862.2
$ws.Range("K138").Value = 5251.1052
$ws.Range("L138").Value = 11586.6
$ws.Range("M138").Value = -111.1052
$ws.Range("N138").Value = -21866.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 825.86365
$ws.Range("I5").Value = 885.3125
$ws.Range("J5").Value = 667.3333
$ws.Range("K5").Value = 885.3125
$ws.Range("L5").Value = 667.3333
$ws.Range("M5").Value = -773.3125
$ws.Range("N5").Value = -891.3333

$ws.Range("H45").Value = 81221.234
$ws.Range("I45").Value = 113209.22
$ws.Range("J45").Value = 9248.25
$ws.Range("K45").Value = 113209.22
$ws.Range("L45").Value = 9248.25
$ws.Range("M45").Value = -112832.22
$ws.Range("N45").Value = -10002.25

$ws.Range("H122").Value = 7410189
$ws.Range("I122").Value = 9664207
$ws.Range("K122").Value = 28992621
$ws.Range("M122").Value = -28990171

$ws.Range("H132").Value = 1965.6364
$ws.Range("I132").Value = 1289.3529
$ws.Range("J132").Value = 2684.1875
$ws.Range("K132").Value = 3868.0587
$ws.Range("L132").Value = 8052.5625
$ws.Range("M132").Value = -1338.0587
$ws.Range("N132").Value = -13112.5625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 825.86365
$ws.Range("I4").Value = 885.3125
$ws.Range("J4").Value = 667.3333
$ws.Range("K4").Value = 885.3125
$ws.Range("L4").Value = 667.3333
$ws.Range("M4").Value = -770.3125
$ws.Range("N4").Value = -897.3333

$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H99").Value = 2791.5715
$ws.Range("I99").Value = 2007.8
$ws.Range("J99").Value = 4751
$ws.Range("K99").Value = 2007.8
$ws.Range("L99").Value = 4751
$ws.Range("M99").Value = -509.8
$ws.Range("N99").Value = -7747

$ws.Range("H105").Value = 1619.4
$ws.Range("I105").Value = 1655.7142
$ws.Range("J105").Value = 1111
$ws.Range("K105").Value = 1655.7142
$ws.Range("L105").Value = 1111
$ws.Range("M105").Value = 91.28580000000011
$ws.Range("N105").Value = -4605

$ws.Range("H130").Value = 53853.332
$ws.Range("J130").Value = 53853.332
$ws.Range("L130").Value = 53853.332
$ws.Range("N130").Value = -63893.332

$ws.Range("H134").Value = 2592.9424
$ws.Range("I134").Value = 1377.9459
$ws.Range("J134").Value = 5589.933
$ws.Range("K134").Value = 4133.8377
$ws.Range("L134").Value = 16769.799
$ws.Range("M134").Value = -1598.8377
$ws.Range("N134").Value = -21839.799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13948.916
$ws.Range("I31").Value = 1523.6
$ws.Range("J31").Value = 15410.718
$ws.Range("K31").Value = 1523.6
$ws.Range("L31").Value = 15410.718
$ws.Range("M31").Value = -1228.6
$ws.Range("N31").Value = -16000.718

$ws.Range("H34").Value = 13948.916
$ws.Range("I34").Value = 1523.6
$ws.Range("J34").Value = 15410.718
$ws.Range("K34").Value = 1523.6
$ws.Range("L34").Value = 15410.718
$ws.Range("M34").Value = -1321.6
$ws.Range("N34").Value = -15814.718

$ws.Range("H99").Value = 4180.8184
$ws.Range("I99").Value = 3200
$ws.Range("J99").Value = 4998.1665
$ws.Range("K99").Value = 3200
$ws.Range("L99").Value = 4998.1665
$ws.Range("M99").Value = -1702
$ws.Range("N99").Value = -7994.1665

$ws.Range("H126").Value = 4180.8184
$ws.Range("I126").Value = 3200
$ws.Range("J126").Value = 4998.1665
$ws.Range("K126").Value = 9600
$ws.Range("L126").Value = 14994.4995
$ws.Range("M126").Value = -7130
$ws.Range("N126").Value = -19934.4995

$ws.Range("H127").Value = 150000
$ws.Range("J127").Value = 150000
$ws.Range("L127").Value = 150000
$ws.Range("N127").Value = -159920

$ws.Range("H132").Value = 40559.625
$ws.Range("I132").Value = 3596.8572
$ws.Range("K132").Value = 10790.5716
$ws.Range("M132").Value = -8260.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1287.9
$ws.Range("I3").Value = 984.875
$ws.Range("K3").Value = 2954.625
$ws.Range("M3").Value = -2842.625

$ws.Range("H56").Value = 19236290
$ws.Range("I56").Value = 19236290
$ws.Range("K56").Value = 19236290
$ws.Range("M56").Value = -19235760

$ws.Range("H68").Value = 437.66666
$ws.Range("I68").Value = 441.2
$ws.Range("K68").Value = 1323.6
$ws.Range("M68").Value = -512.5999999999999

$ws.Range("H71").Value = 437.66666
$ws.Range("I71").Value = 441.2
$ws.Range("K71").Value = 3970.8
$ws.Range("M71").Value = 85.20000000000027

$ws.Range("H110").Value = 14144.77
$ws.Range("J110").Value = 17688.7
$ws.Range("L110").Value = 53066.10000000001
$ws.Range("N110").Value = -61246.10000000001

$ws.Range("H140").Value = 2246.5908
$ws.Range("I140").Value = 2071.3
$ws.Range("K140").Value = 6213.900000000001
$ws.Range("M140").Value = -1033.900000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 25262.273
$ws.Range("I102").Value = 2020.5312
$ws.Range("J102").Value = 87240.25
$ws.Range("K102").Value = 2020.5312
$ws.Range("L102").Value = 87240.25
$ws.Range("M102").Value = -398.5311999999999
$ws.Range("N102").Value = -90484.25

$ws.Range("H122").Value = 205680.66
$ws.Range("J122").Value = 4844.4707
$ws.Range("L122").Value = 14533.4121
$ws.Range("N122").Value = -19433.4121

$ws.Range("H132").Value = 3668.739
$ws.Range("I132").Value = 3095.1428
$ws.Range("J132").Value = 4561
$ws.Range("K132").Value = 9285.428400000001
$ws.Range("L132").Value = 13683
$ws.Range("M132").Value = -6755.428400000001
$ws.Range("N132").Value = -18743

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 69000.07000000001
$ws.Range("I22").Value = 100851.2
$ws.Range("J22").Value = 5297.8
$ws.Range("K22").Value = 100851.2
$ws.Range("L22").Value = 5297.8
$ws.Range("M22").Value = -100556.2
$ws.Range("N22").Value = -5887.8

$ws.Range("H27").Value = 69000.07000000001
$ws.Range("I27").Value = 100851.2
$ws.Range("J27").Value = 5297.8
$ws.Range("K27").Value = 100851.2
$ws.Range("L27").Value = 5297.8
$ws.Range("M27").Value = -100744.2
$ws.Range("N27").Value = -5511.8

$ws.Range("H48").Value = 26333.334
$ws.Range("I48").Value = 9000
$ws.Range("K48").Value = 9000
$ws.Range("M48").Value = -8339

$ws.Range("H61").Value = 2893.3333
$ws.Range("I61").Value = 2893.3333
$ws.Range("K61").Value = 2893.3333
$ws.Range("M61").Value = -2691.3333

$ws.Range("H100").Value = 3379.25
$ws.Range("I100").Value = 3147.7144
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 3147.7144
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -2606.7144
$ws.Range("N100").Value = -6082

$ws.Range("H113").Value = 2893.3333
$ws.Range("I113").Value = 2893.3333
$ws.Range("K113").Value = 2893.3333
$ws.Range("M113").Value = -723.3332999999998

$ws.Range("H122").Value = 5934.684
$ws.Range("J122").Value = 7327.1
$ws.Range("L122").Value = 21981.3
$ws.Range("N122").Value = -26881.3

$ws.Range("H129").Value = 120000
$ws.Range("J129").Value = 120000
$ws.Range("L129").Value = 120000
$ws.Range("N129").Value = -130000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 49000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 49000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 49000
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -59000

$ws.Range("H132").Value = 21881.32
$ws.Range("I132").Value = 4761.086
$ws.Range("J132").Value = 61828.535
$ws.Range("K132").Value = 14283.258
$ws.Range("L132").Value = 185485.605
$ws.Range("M132").Value = -11753.258
$ws.Range("N132").Value = -190545.605

$ws.Range("H136").Value = 3319.2173
$ws.Range("I136").Value = 2520.2354
$ws.Range("J136").Value = 5583
$ws.Range("K136").Value = 7560.706200000001
$ws.Range("L136").Value = 16749
$ws.Range("M136").Value = -5010.706200000001
$ws.Range("N136").Value = -21849
